$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# New row 29: WIN / CONGRADUALATION! YOU WIN!!! / CHÚC MỪNG BẠN ĐÃ CHIẾN THẮNG
$ws.Range("A29").Value = "WIN"
$ws.Range("B29").Value = "CONGRADUALATION! YOU WIN!!!"
$ws.Range("C29").Value = "CHÚC MỪNG BẠN ĐÃ CHIẾN THẮNG"

# Match the saved view state: selection on C29, scrolled so row 7 is at the top
$ws.Range("C29").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
